$wb = $excel.ActiveWorkbook

# --- Teachers sheet: insert two new teacher rows ("Achille Cochet" and
# "Henry D'Aboville", both assigned to UNI011-UNI021) right above the
# "Christiane Brunelle" row, pushing the remaining rows down. ---
$teachers = $wb.Worksheets.Item("Teachers")

$teachers.Rows.Item(5).Insert()
$teachers.Rows.Item(5).Insert()

$teachers.Range("A5").Value = 3
$teachers.Range("B5").Value = "Achille"
$teachers.Range("C5").Value = "Cochet"
$teachers.Range("D5").Value = "UNI011-UNI021"

$teachers.Range("A6").Value = 4
$teachers.Range("B6").Value = "Henry"
$teachers.Range("C6").Value = "D'Aboville"
$teachers.Range("D6").Value = "UNI011-UNI021"

# Renumber the "Id" column (A) for the rows that got pushed down.
$teachers.Range("A7").Value = 5
$teachers.Range("A8").Value = 6
$teachers.Range("A9").Value = 7
$teachers.Range("A10").Value = 8
$teachers.Range("A11").Value = 9

$teachers.Rows.Item(3).Select()

# --- Selections on other sheets (match the final saved cursor positions) ---
$promotions = $wb.Worksheets.Item("Promotions")
$promotions.Range("A4").Select()

$rooms = $wb.Worksheets.Item("Rooms")
$rooms.Activate()
$rooms.Range("H8").Select()
